$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Leading apostrophe forces text entry so values match the source feed's
# string formatting (e.g. "308.44", "1.65%") instead of Excel-coerced numbers.
$ws.Range("D2").Value = "'308.44"
$ws.Range("E2").Value = "'1.65%"
$ws.Range("D3").Value = "'39.12"
$ws.Range("E3").Value = "'9.55%"
$ws.Range("D4").Value = "'5.117"
$ws.Range("E4").Value = "'1.58%"
$ws.Range("D5").Value = "'0.08162"
$ws.Range("E5").Value = "'3.40%"
$ws.Range("D6").Value = "'2.048"
$ws.Range("E6").Value = "'12.06%"
$ws.Range("D7").Value = "'4.177"
$ws.Range("E7").Value = "'1.70%"
$ws.Range("D8").Value = "'7.900"
$ws.Range("E8").Value = "'1.47%"
$ws.Range("D9").Value = "'0.9297"
$ws.Range("E9").Value = "'1.02%"
$ws.Range("D10").Value = "'0.1403"
$ws.Range("E10").Value = "'4.11%"
$ws.Range("D11").Value = "'0.1945"
$ws.Range("E11").Value = "'2.65%"
$ws.Range("D12").Value = "'0.09205"
$ws.Range("E12").Value = "'1.23%"
$ws.Range("D13").Value = "'0.03487"
$ws.Range("E13").Value = "'0.52%"
$ws.Range("D14").Value = "'0.09844"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("D15").Value = "'0.001408"
$ws.Range("E15").Value = "'-1.01%"
$ws.Range("D16").Value = "'0.005864"
$ws.Range("E16").Value = "'-3.32%"
$ws.Range("D17").Value = "'3.946"
$ws.Range("E17").Value = "'6.02%"
$ws.Range("D18").Value = "'3.445"
$ws.Range("E18").Value = "'3.01%"
$ws.Range("D19").Value = "'0.3454"
$ws.Range("E19").Value = "'0.48%"
$ws.Range("D20").Value = "'0.1302"
$ws.Range("E20").Value = "'-0.58%"
$ws.Range("D21").Value = "'4.796"
$ws.Range("E21").Value = "'-7.50%"
$ws.Range("D23").Value = "'0.04467"
$ws.Range("E23").Value = "'1.20%"
$ws.Range("D24").Value = "'0.001242"
$ws.Range("E24").Value = "'0.46%"
$ws.Range("E25").Value = "'-9.56%"
$ws.Range("E27").Value = "'-0.16%"
$ws.Range("D39").Value = "'0.02113"
$ws.Range("E39").Value = "'9.04%"
$ws.Range("D40").Value = "'0.05167"
$ws.Range("E40").Value = "'1.79%"
$ws.Range("D41").Value = "'0.007467"
$ws.Range("E41").Value = "'-2.04%"
$ws.Range("D42").Value = "'0.01013"
$ws.Range("E42").Value = "'-0.38%"
$ws.Range("E43").Value = "'1.87%"
$ws.Range("D44").Value = "'0.002129"
$ws.Range("E44").Value = "'-1.55%"
$ws.Range("D45").Value = "'0.009669"
$ws.Range("E45").Value = "'-5.10%"
$ws.Range("D46").Value = "'0.00006313"
$ws.Range("E46").Value = "'1.85%"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E48").Value = "'-0.63%"
$ws.Range("D49").Value = "'0.001601"
$ws.Range("E49").Value = "'-3.60%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.05%"
